$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("Лист2")

# 1) Fix typo: "Карачаево Черкесская Республика" -> "Карачаево-Черкесская Республика"
$karachay = $ws2.Range("A1:B100").Find("Карачаево Черкесская Республика")
if ($karachay -ne $null) {
    $karachay.Value = "Карачаево-Черкесская Республика"
}

# 2) Insert new row "Ненецкий автономный округ" right before "Ленинградская область"
$lenRow = $ws2.Range("A1:B100").Find("Ленинградская область")
$lenRowIndex = $lenRow.Row
$ws2.Rows.Item($lenRowIndex).Insert()
$ws2.Range("A" + $lenRowIndex).Value = $ws2.Range("A" + ($lenRowIndex + 1)).Value()
$ws2.Range("B" + $lenRowIndex).Value = "Ненецкий автономный округ"

# 3) Insert new row "Республика Калмыкия" right before "Республика Дагестан"
$dagRow = $ws2.Range("A1:B100").Find("Республика Дагестан")
$dagRowIndex = $dagRow.Row
$ws2.Rows.Item($dagRowIndex).Insert()
$ws2.Range("A" + $dagRowIndex).Value = $ws2.Range("A" + ($dagRowIndex + 1)).Value()
$ws2.Range("B" + $dagRowIndex).Value = "Республика Калмыкия"

# 4) Update view state: Лист2 becomes the active/selected sheet
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 53
$ws2.Range("B75").Select()

# 5) Update view state on Лист1 (no longer the active tab, selection moves to B32)
$ws1.Range("B32").Select()
